$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Start of Project"
$ws.Range("B3").Value = "Standard A"
$ws.Range("B4").Value = "Inverted Cosmonauts"
$ws.Range("B5").Value = "Start of Construction/build"
